$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates are serial numbers, 2021-08-24 .. 2021-09-01)
$rows = @(
    @{ Row = 358; A = 44432; B = 0; C = 3; D = 112.4016485575122 },
    @{ Row = 359; A = 44433; B = 0; C = 3; D = 112.4016485575122 },
    @{ Row = 360; A = 44434; B = 0; C = 3; D = 112.4016485575122 },
    @{ Row = 361; A = 44435; B = 0; C = 2; D = 74.93443237167479 },
    @{ Row = 362; A = 44436; B = 0; C = 2; D = 74.93443237167479 },
    @{ Row = 363; A = 44437; B = 0; C = 1; D = 37.46721618583739 },
    @{ Row = 364; A = 44438; B = 0; C = 0; D = 0 },
    @{ Row = 365; A = 44439; B = 0; C = 0; D = 0 },
    @{ Row = 366; A = 44440; B = 0; C = 0; D = 0 }
)

# Last existing data row (357) carries the date cell style (border, bold,
# centered, custom date numFmt) that new date cells in column A must copy.
$styleSource = $ws.Cells.Item(357, 1)
$styleSource.Copy()

foreach ($r in $rows) {
    $cellA = $ws.Cells.Item($r.Row, 1)
    $cellA.PasteSpecial(-4122)
    $cellA.Value = $r.A

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
